# Rebuild the "Locher CAZy clan" conversion table with the new, expanded
# layout: a "Structural subclass Locher" column (A, merged per group), the
# existing "CAZy clan" / "CAZy families" columns (B / C), and two new
# columns for inverting/retaining mechanism (D) and donor substrate (E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- clear the cell that is no longer part of the header row -------------
$ws.Range("D1").Value = ""

# --- header row ------------------------------------------------------------
$ws.Range("A1").Value = "Structural subclass Locher"
$ws.Range("B1").Value = "CAZy clan"
$ws.Range("C1").Value = "CAZy families"
$ws.Range("E1").Value = "Activator"

# --- GT-CB (10 conserved) block: rows 2-5 ----------------------------------
$ws.Range("A2").Value = "GT-CB (10 conserved)"
$ws.Range("B2").Value = "GT-CB1"
$ws.Range("C2").Value = "X605, X607, X609, X613, X614, X615"
$ws.Range("D2").Value = "Inverting"
$ws.Range("E2").Value = "Lipid-PP-oligosaccharide"

$ws.Range("B3").Value = "GT-CB2"
$ws.Range("C3").Value = "X586, X606, X608, X610, X611, X612"
$ws.Range("D3").Value = "Retaining"
$ws.Range("E3").Value = "Lipid-PP-oligosaccharide"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "X571"
$ws.Range("D4").Value = "Inverting"
$ws.Range("E4").Value = "Lipid-PP-oligosaccharide"

$ws.Range("B5").Value = "-"
$ws.Range("C5").Value = "X617"
$ws.Range("D5").Value = "Inverting"
$ws.Range("E5").Value = "Lipid-PP-oligosaccharide"

# --- GT-CA (7 conserved) block: rows 6-10 ----------------------------------
$ws.Range("A6").Value = "GT-CA (7 conserved)"
$ws.Range("B6").Value = "-"
$ws.Range("C6").Value = "GT66"
$ws.Range("D6").Value = "Inverting"
$ws.Range("E6").Value = "Lipid-PP-oligosaccharide"

$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "GT83"
$ws.Range("D7").Value = "Inverting"
$ws.Range("E7").Value = "Lipid-P"

$ws.Range("B8").Value = "-"
$ws.Range("C8").Value = "GT39"
$ws.Range("D8").Value = "Inverting"
$ws.Range("E8").Value = "Lipid-P"

$ws.Range("B9").Value = "-"
$ws.Range("C9").Value = "GT57"
$ws.Range("D9").Value = "Inverting"
$ws.Range("E9").Value = "Lipid-P"

$ws.Range("B10").Value = "-"
$ws.Range("C10").Value = "GT53"
$ws.Range("D10").Value = "Inverting"

# --- ungrouped rows 11-14 (no Locher subclass merge) -----------------------
$ws.Range("A11").Value = "-"
$ws.Range("B11").Value = "-"
$ws.Range("C11").Value = "GT22"
$ws.Range("D11").Value = "Inverting"
$ws.Range("E11").Value = "Lipid-P"

$ws.Range("A12").Value = "-"
$ws.Range("B12").Value = "-"
$ws.Range("C12").Value = "GT50"
$ws.Range("D12").Value = "Inverting"
$ws.Range("E12").Value = "Lipid-P"

$ws.Range("A13").Value = "-"
$ws.Range("B13").Value = "-"
$ws.Range("C13").Value = "GT58"
$ws.Range("D13").Value = "Inverting"
$ws.Range("E13").Value = "Lipid-P"

$ws.Range("A14").Value = "-"
$ws.Range("B14").Value = "-"
$ws.Range("C14").Value = "GT59"
$ws.Range("D14").Value = "Inverting"
$ws.Range("E14").Value = "Lipid-P"

# --- merge + center the "Structural subclass Locher" group headers --------
# Setting alignment on the single anchor cell before merging (rather than on
# the whole multi-cell range) propagates the centred style to every cell in
# the merged block without Excel fabricating extra intermediate style
# records for each property assignment on the range. The second group reuses
# the first group's format (copy/paste-special of formats only) instead of
# re-deriving it through HorizontalAlignment/VerticalAlignment again, so
# both groups end up sharing one cell style.
$ws.Cells.Item(2, 1).HorizontalAlignment = -4108
$ws.Cells.Item(2, 1).VerticalAlignment = -4108
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial(-4122)

$ws.Range("A2:A5").Merge()
$ws.Range("A6:A10").Merge()

# --- column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.4987
$ws.Columns.Item(3).ColumnWidth = 30.8307
$ws.Columns.Item(5).ColumnWidth = 21.1667

# --- selection, matching the saved workbook's cursor position -------------
$ws.Range("C19").Select()

Write-Output "conversion-locher table rebuilt"
